# Horarios actualizados Linea 141 - 610
# Apply the scraped-data refresh: updated timestamps / totals, re-sorted
# rows with swapped/changed values, and newly appended rows at the end
# of each of the three schedule sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 19:11:44"
$ws.Cells.Item(3, 1).Value = "Total filas: 326"
$ws.Cells.Item(45, 1).Value = "07:56:02"
$ws.Cells.Item(45, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(45, 4).Value = 49
$ws.Cells.Item(46, 1).Value = "08:45:31"
$ws.Cells.Item(46, 3).Value = "215C_EL PATO"
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(66, 1).Value = "08:28:52"
$ws.Cells.Item(66, 3).Value = "10_OLMOS"
$ws.Cells.Item(66, 4).Value = 60
$ws.Cells.Item(67, 1).Value = "08:11:18"
$ws.Cells.Item(67, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(67, 4).Value = 77
$ws.Cells.Item(86, 1).Value = "08:38:24"
$ws.Cells.Item(86, 3).Value = "15_ABASTO"
$ws.Cells.Item(86, 4).Value = 111
$ws.Cells.Item(87, 1).Value = "10:04:30"
$ws.Cells.Item(87, 3).Value = "14_ABASTO"
$ws.Cells.Item(87, 4).Value = 25
$ws.Cells.Item(106, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(107, 3).Value = "10_OLMOS"
$ws.Cells.Item(124, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(125, 3).Value = "15_ABASTO"
$ws.Cells.Item(137, 1).Value = "11:46:32"
$ws.Cells.Item(137, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(137, 4).Value = 48
$ws.Cells.Item(138, 1).Value = "10:36:50"
$ws.Cells.Item(138, 3).Value = "15_ABASTO"
$ws.Cells.Item(138, 4).Value = 118
$ws.Cells.Item(147, 1).Value = "11:13:15"
$ws.Cells.Item(147, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(147, 4).Value = 110
$ws.Cells.Item(148, 1).Value = "11:33:52"
$ws.Cells.Item(148, 3).Value = "215C_EL PATO"
$ws.Cells.Item(148, 4).Value = 90
$ws.Cells.Item(158, 1).Value = "12:11:21"
$ws.Cells.Item(158, 3).Value = "14_ABASTO"
$ws.Cells.Item(158, 4).Value = 81
$ws.Cells.Item(159, 1).Value = "11:53:44"
$ws.Cells.Item(159, 3).Value = "215A_EL PATO"
$ws.Cells.Item(159, 4).Value = 99
$ws.Cells.Item(178, 1).Value = "12:33:02"
$ws.Cells.Item(178, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(178, 4).Value = 104
$ws.Cells.Item(179, 1).Value = "12:53:26"
$ws.Cells.Item(179, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(179, 4).Value = 84
$ws.Cells.Item(258, 1).Value = "16:37:37"
$ws.Cells.Item(258, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(258, 4).Value = 76
$ws.Cells.Item(259, 1).Value = "16:51:51"
$ws.Cells.Item(259, 3).Value = "10_OLMOS"
$ws.Cells.Item(259, 4).Value = 62
$ws.Cells.Item(269, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(270, 3).Value = "15_ABASTO"
$ws.Cells.Item(301, 1).Value = "17:55:25"
$ws.Cells.Item(301, 3).Value = "17_ROMERO"
$ws.Cells.Item(301, 4).Value = 81
$ws.Cells.Item(302, 1).Value = "17:35:41"
$ws.Cells.Item(302, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(302, 4).Value = 101
$ws.Cells.Item(323, 1).Value = "19:11:44"
$ws.Cells.Item(323, 2).Value = "20:13"
$ws.Cells.Item(323, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(323, 4).Value = 62
$ws.Cells.Item(324, 1).Value = "18:30:48"
$ws.Cells.Item(324, 2).Value = "20:21"
$ws.Cells.Item(324, 4).Value = 111
$ws.Cells.Item(325, 1).Value = "18:52:29"
$ws.Cells.Item(325, 2).Value = "20:22"
$ws.Cells.Item(325, 3).Value = "15_ABASTO"
$ws.Cells.Item(325, 4).Value = 90
$ws.Cells.Item(326, 1).Value = "18:44:45"
$ws.Cells.Item(326, 2).Value = "20:30"
$ws.Cells.Item(326, 3).Value = "10_OLMOS"
$ws.Cells.Item(326, 4).Value = 106
$ws.Cells.Item(327, 1).Value = "19:11:44"
$ws.Cells.Item(327, 2).Value = "20:41"
$ws.Cells.Item(327, 3).Value = "17_ROMERO"
$ws.Cells.Item(327, 4).Value = 90
$ws.Cells.Item(328, 1).Value = "18:52:29"
$ws.Cells.Item(328, 2).Value = "20:42"
$ws.Cells.Item(328, 3).Value = "17_ROMERO"
$ws.Cells.Item(328, 4).Value = 110
$ws.Cells.Item(328, 5).Value = "LP1912"
$ws.Cells.Item(329, 1).Value = "18:52:29"
$ws.Cells.Item(329, 2).Value = "20:47"
$ws.Cells.Item(329, 3).Value = "215B_EL PATO"
$ws.Cells.Item(329, 4).Value = 115
$ws.Cells.Item(329, 5).Value = "LP1912"
$ws.Cells.Item(330, 1).Value = "19:11:44"
$ws.Cells.Item(330, 2).Value = "20:56"
$ws.Cells.Item(330, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(330, 4).Value = 105
$ws.Cells.Item(330, 5).Value = "LP1912"
$ws.Cells.Item(331, 1).Value = "19:11:44"
$ws.Cells.Item(331, 2).Value = "21:06"
$ws.Cells.Item(331, 3).Value = "10_OLMOS"
$ws.Cells.Item(331, 4).Value = 115
$ws.Cells.Item(331, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 19:11:44"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 19:11:44"
$ws.Cells.Item(3, 1).Value = "Total filas: 46"
$ws.Cells.Item(51, 1).Value = "19:11:44"
$ws.Cells.Item(51, 2).Value = "20:51"
$ws.Cells.Item(51, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(51, 4).Value = 100
$ws.Cells.Item(51, 5).Value = "L6203"
